$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2   = 1
    29  = 1
    39  = 0
    49  = 0
    60  = 0
    90  = 0
    96  = 0
    124 = 0
    145 = 0
    201 = 0
    202 = 0
    208 = 1
    222 = 1
    231 = 1
    254 = 0
    278 = 0
    280 = 1
    294 = 1
    376 = 1
    389 = 0
    393 = 0
    456 = 1
    481 = 0
    488 = 1
    503 = 0
    523 = 0
    528 = 1
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 1).Value = $changes[$row]
}
